# LOM3202.xlsx edit
# 1) Remove row 13 (the row that only had B/C = "5982760 - Carlos Alberto Baldan",
#    with no label in column A). This shifts every following row up by one and
#    shrinks the used range from A1:C25 to A1:C24.
# 2) Update the text of several (now-shifted) B/C cells so the sheet ends up
#    matching the published content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 13 entirely (shifts rows 14-25 up to 13-24).
$ws.Rows(13).Delete()

# Column B and C share the same text per row, so update both together.
$ws.Range("B10:C10").Value = "5982760 - Carlos Alberto Baldan"
$ws.Range("B13:C13").Value = "Semestral"
$ws.Range("B15:C15").Value = "01/01/2015"
$ws.Range("B18:C18").Value = "5982760 - Carlos Alberto Baldan"
$ws.Range("B19:C19").Value = "Aulas expositivas teóricas, aulas de exercícios e estudos de casos por simulação computacional."
$ws.Range("B20:C20").Value = "Duas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + 2P2)/3"
$ws.Range("B21:C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

Write-Host "Done"
